# Insert a new weekly price record for "Perejil" (Vega Central Mapocho de
# Santiago) ahead of the existing row 295, shifting all subsequent rows
# down by one (old row 295 becomes 296, ..., old row 313 becomes 314).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 295; Excel shifts rows 295:313 down
# to 296:314 and copies formatting (including the date style on column D)
# from the row above, matching the workbook's existing pattern.
$ws.Rows("295:295").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A295").Value = 9
$ws.Range("B295").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C295").Value = "Metropolitana"
$ws.Range("D295").Value = 44610
$ws.Range("E295").Value = 13
$ws.Range("F295").Value = 100112044
$ws.Range("G295").Value = "Perejil"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 61
$ws.Range("K295").Value = 18000
$ws.Range("L295").Value = 19000
$ws.Range("M295").Value = 18508
$ws.Range("N295").Value = "$/docena de atados"
$ws.Range("O295").Value = "Región Metropolitana"
$ws.Range("P295").Value = 6169
$ws.Range("Q295").Value = 3
$ws.Range("R295").Value = "Hortaliza"
